# Update the timestamp in column A for the newly appended rows (2-7)
# on the "ランサーズ" sheet to reflect the latest fetch time.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-06 18:26:44"

for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
